# Change to total stress and strain formulation instead of incremental form
# as in the original formulation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core formula changes (rows 3-18 unless noted) -------------------------

# AS2 was a literal 0; it now computes the (former AR3-style) total-strain
# expression using the NEXT row's inputs.
$ws.Range("AS2").Formula = "=(1+2*AM3)*AK3*(1-AP3/3)/(2*AN3*AO3*AM3)"

for ($r = 3; $r -le 18; $r++) {
    $rNext = $r + 1

    # AI: ratio now uses total (AD/AE) instead of incremental (AG/AH) strains
    $ws.Cells.Item($r, 35).Formula = "=-AD$r/AE$r"

    # AP: ratio^0.25 computed directly from AJ (total), instead of the
    # incremental secant slope between consecutive rows
    $ws.Cells.Item($r, 42).Formula = "=AJ$r^0.25"

    # AS: total strain term uses the NEXT row's own values instead of
    # accumulating row-to-row increments
    $ws.Cells.Item($r, 45).Formula = "=(1+2*AM$rNext)*AK$rNext*(1-AP$rNext/3)/(2*AN$rNext*AO$rNext*AM$rNext)"

    # AU: total strain term uses THIS row's own values instead of
    # accumulating row-to-row increments
    $ws.Cells.Item($r, 47).Formula = "=(1-AM$r)*(AK$r)*(1-AP$r/3)/(3*AN$r*AO$r*AM$r)"
}

# --- Highlight the recomputed cells with the existing yellow fill style ----
# (reuses the workbook's existing "s=1" style; mirrors the diff which shows
# these ranges switching to the yellow-fill style without adding new styles)
$ws.Range("AS2").Interior.Color = 65535
for ($r = 3; $r -le 18; $r++) {
    $ws.Range("AI$r").Interior.Color = 65535
    $ws.Range("AO$r").Interior.Color = 65535
    $ws.Range("AP$r").Interior.Color = 65535
    $ws.Range("AS$r").Interior.Color = 65535
    $ws.Range("AU$r").Interior.Color = 65535
}

# --- Cosmetic workbook / view changes ---------------------------------------

$ws.Application.ActiveWindow.ScrollColumn = 29  ## AC
$excel.ActiveWindow.Zoom = 55
$ws.Range("AI1:AI1048576").Select()

Write-Host "Edit complete"
